# "Loading in newest US data" -- update the TTLE workbook's logit-exponent
# notes and values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# About sheet: refresh the calibration notes.
# ---------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

# New note about the non-road logit exponent (set before the "onroad" note
# below so the new shared strings are introduced in the same order as the
# authored workbook).
$wsAbout.Range("A13").Value = "For non-road we use -0.1 because of costs preventing the model from solving."

# Replaces the old PNNL GCAM calibration note.
$wsAbout.Range("A12").Value = "We use calibrated values in onroad sectors."

# "Modified Logit" -> "Unmodified Logit".
$wsAbout.Range("A15").Value = 'For more on this, see the "Unmodified Logit" equation description at:'

$wsAbout.Range("A16").Select() | Out-Null

# ---------------------------------------------------------------------
# TTLE sheet: load in the newest US logit exponent values.
# ---------------------------------------------------------------------
$wsTTLE = $wb.Worksheets.Item("TTLE")
$wsTTLE.Activate() | Out-Null

$wsTTLE.Range("B2").Value = -80    # LDVs, Passenger

$wsTTLE.Range("B4").Value = -0.1   # aircraft, Passenger
$wsTTLE.Range("C4").Value = -0.1   # aircraft, Freight

$wsTTLE.Range("B5").Value = -0.15  # rail, Passenger
$wsTTLE.Range("C5").Value = -0.1   # rail, Freight

$wsTTLE.Range("B6").Value = -0.1   # ships, Passenger
$wsTTLE.Range("C6").Value = -0.1   # ships, Freight

$wsTTLE.Range("G3").Select() | Out-Null
